# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (fund-level holdings) right after "2021-Q4"
#    and before "总计", built from a copy of "2021-Q2" (same column layout,
#    one row short) so that fonts / borders / style indexes line up with the
#    other quarter sheets.
# 2. Insert a new first data row into "总计" summarizing the 2022-Q1 totals,
#    pushing the previous rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the "2022-Q1" worksheet
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q2")
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($null, $afterSheet)

$q1 = $wb.Worksheets.Item("2021-Q2 (2)")
$q1.Name = "2022-Q1"

# The template only has 5 data rows (rows 2-6); we need 6 (rows 2-7).
# Duplicate the formatting of the last data row into the new row 7.
$q1.Range("A6:H6").Copy($q1.Range("A7:H7"))

# --- header row (plain text assignment; keeps the bordered header style
#     that was copied over from the template sheet) ---
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- data rows ---
$q1Data = @(
    @(0, "162102", "金鹰中小盘精选混合",         "4.60", "76.52", "4.45", "0.2047", 3),
    @(1, "001167", "金鹰科技创新股票",           "4.03", "94.55", "5.03", "0.2027", 4),
    @(2, "210009", "金鹰核心资源混合",           "3.86", "94.96", "4.82", "0.1861", 3),
    @(3, "001411", "诺安创新驱动灵活配置混合A",   "3.96", "94.02", "4.36", "0.1727", 10),
    @(4, "004044", "金鹰转型动力灵活配置混合",     "0.72", "93.34", "8.64", "0.0622", 1),
    @(5, "002051", "诺安创新驱动灵活配置混合C",   "1.33", "94.02", "4.36", "0.0580", 10)
)

# The B:G columns hold numeric-looking codes/figures (e.g. "162102",
# "4.60") that must stay text (t="inlineStr"/t="s"), matching the source
# file, instead of being auto-converted to numbers. Temporarily mark the
# whole block as Text, fill in the values, then clear the formatting back
# off again (one NumberFormat/ClearFormats pair for the whole block, so it
# doesn't fragment the style table).
$q1TextRange = $q1.Range("B2:G7")
$q1TextRange.NumberFormat = "@"
$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    $q1.Range("G$r").Value = $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}
$q1TextRange.ClearFormats()

# ---------------------------------------------------------------------
# Step 2: update the "总计" worksheet with a new 2022-Q1 summary row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

$total.Range("A2").Value = 0
$totalTextRange = $total.Range("B2")
$totalTextRange.NumberFormat = "@"
$totalTextRange.Value = "2022-Q1"
$totalTextRange.ClearFormats()
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.89

# The index column (A) is a running counter; renumber the rows that were
# pushed down (they used to start at 0, now they start at 1).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the originally active sheet/tab so this edit doesn't change the
# workbook's selected-sheet state.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q1 sheet and 总计 summary row added"
